# Scheduled market-data refresh for the Bahamut Profits workbook.
# Universalis price snapshots changed since the last run, so the
# cached currentAveragePrice* / LevePrice* / LeveProfit* columns
# are refreshed per leve row below (one sheet per crafting job).

$wb = $excel.ActiveWorkbook

# ALC!15 - "Morning Glass of Ether" (Ether)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 1825.0618
$ws.Range("I15").Value = 1825.0618
$ws.Range("K15").Value = 5475.1854
$ws.Range("M15").Value = -5306.1854

# ALC!88 - "The Grave of Hemlock Groves" (Growth Formula Zeta)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H88").Value = 3969003.8
$ws.Range("J88").Value = 5291671.5
$ws.Range("L88").Value = 5291671.5
$ws.Range("N88").Value = -5292483.5

# ALC!91 - "Dappling the Highlands (L)" (Growth Formula Zeta)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H91").Value = 3969003.8
$ws.Range("J91").Value = 5291671.5
$ws.Range("L91").Value = 5291671.5
$ws.Range("N91").Value = -5294479.5

# ALC!113 - "Amaro Kart" (Starch Glue)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H113").Value = 3341.1
$ws.Range("I113").Value = 7502.5
$ws.Range("J113").Value = 2300.75
$ws.Range("K113").Value = 7502.5
$ws.Range("L113").Value = 2300.75
$ws.Range("M113").Value = -4248.5
$ws.Range("N113").Value = -8808.75

# ALC!132 - "Fast-forwarding Flora" (Growth Formula Lambda)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 2878.4138
$ws.Range("J132").Value = 3987.5
$ws.Range("L132").Value = 11962.5
$ws.Range("N132").Value = -17022.5

# ALC!137 - "Cutting Edge of Culinary Quality" (Magnesia Whetstone)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 2267.543
$ws.Range("I137").Value = 2345.6155
$ws.Range("J137").Value = 2221.4092
$ws.Range("K137").Value = 7036.8465
$ws.Range("L137").Value = 6664.2276
$ws.Range("M137").Value = -4486.8465
$ws.Range("N137").Value = -11764.2276

# ALC!138 - "All-night Crafting" (Cunning Craftsman's Tisane)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 2462.13
$ws.Range("I138").Value = 1476.6285
$ws.Range("J138").Value = 2992.7847
$ws.Range("K138").Value = 4429.8855
$ws.Range("L138").Value = 8978.3541
$ws.Range("M138").Value = 710.1144999999997
$ws.Range("N138").Value = -19258.3541

# ARM!2 - "Ain't Got No Ingots" (Bronze Ingot)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 129446.5
$ws.Range("I2").Value = 253330.5
$ws.Range("J2").Value = 5562.5
$ws.Range("K2").Value = 253330.5
$ws.Range("L2").Value = 5562.5
$ws.Range("M2").Value = -253217.5
$ws.Range("N2").Value = -5788.5

# ARM!45 - "Hollow Hallmarks" (Mythril Ingot)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2953.25
$ws.Range("I45").Value = 1568.5
$ws.Range("J45").Value = 4338
$ws.Range("K45").Value = 1568.5
$ws.Range("L45").Value = 4338
$ws.Range("M45").Value = -1191.5
$ws.Range("N45").Value = -5092

# ARM!63 - "Rivets Run through It" (Mythrite Rivets)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 2363.125
$ws.Range("I63").Value = 2415
$ws.Range("J63").Value = 2000
$ws.Range("K63").Value = 2415
$ws.Range("L63").Value = 2000
$ws.Range("M63").Value = -1729
$ws.Range("N63").Value = -3372

# ARM!66 - "A Riveting Revival (L)" (Mythrite Rivets)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H66").Value = 2363.125
$ws.Range("I66").Value = 2415
$ws.Range("J66").Value = 2000
$ws.Range("K66").Value = 12075
$ws.Range("L66").Value = 10000
$ws.Range("M66").Value = -8643
$ws.Range("N66").Value = -16864

# ARM!74 - "As the Bolt Flies" (Titanium Nugget)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 733.8246
$ws.Range("I74").Value = 676.8788
$ws.Range("J74").Value = 812.125
$ws.Range("K74").Value = 676.8788
$ws.Range("L74").Value = 812.125
$ws.Range("M74").Value = 197.1212
$ws.Range("N74").Value = -2560.125

# ARM!77 - "Heavy Metal Banned (L)" (Titanium Nugget)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 733.8246
$ws.Range("I77").Value = 676.8788
$ws.Range("J77").Value = 812.125
$ws.Range("K77").Value = 3384.394
$ws.Range("L77").Value = 4060.625
$ws.Range("M77").Value = 983.6060000000002
$ws.Range("N77").Value = -12796.625

# ARM!97 - "Ore for Me" (High Steel Ingot)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 688.3333
$ws.Range("I97").Value = 688.3333
$ws.Range("K97").Value = 688.3333
$ws.Range("M97").Value = -192.3333

# ARM!116 - "No Scope" (Titanbronze Ingot)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H116").Value = 129446.5
$ws.Range("I116").Value = 253330.5
$ws.Range("J116").Value = 5562.5
$ws.Range("K116").Value = 253330.5
$ws.Range("L116").Value = 5562.5
$ws.Range("M116").Value = -251036.5
$ws.Range("N116").Value = -10150.5

# ARM!122 - "Haste for High Durium" (High Durium Nugget)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 1653
$ws.Range("I122").Value = 1406
$ws.Range("J122").Value = 1900
$ws.Range("K122").Value = 4218
$ws.Range("L122").Value = 5700
$ws.Range("M122").Value = -1768
$ws.Range("N122").Value = -10600

# BSM!3 - "Hells Bells" (Bronze Ingot)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 129446.5
$ws.Range("I3").Value = 253330.5
$ws.Range("J3").Value = 5562.5
$ws.Range("K3").Value = 253330.5
$ws.Range("L3").Value = 5562.5
$ws.Range("M3").Value = -253216.5
$ws.Range("N3").Value = -5790.5

# CRP!31 - "Wall Not Found" (Walnut Lumber)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1643.9
$ws.Range("I31").Value = 835.2
$ws.Range("J31").Value = 2221.543
$ws.Range("K31").Value = 835.2
$ws.Range("L31").Value = 2221.543
$ws.Range("M31").Value = -540.2
$ws.Range("N31").Value = -2811.543

# CRP!34 - "Armoires of the Rich and Famous" (Walnut Lumber)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 1643.9
$ws.Range("I34").Value = 835.2
$ws.Range("J34").Value = 2221.543
$ws.Range("K34").Value = 835.2
$ws.Range("L34").Value = 2221.543
$ws.Range("M34").Value = -633.2
$ws.Range("N34").Value = -2625.543

# CUL!68 - "Such a Butter Face" (Fermented Butter)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 775.75
$ws.Range("I68").Value = 524.9524
$ws.Range("J68").Value = 1202.7838
$ws.Range("K68").Value = 1574.8572
$ws.Range("L68").Value = 3608.3514
$ws.Range("M68").Value = -763.8571999999999
$ws.Range("N68").Value = -5230.3514

# CUL!71 - "No Margarine of Error (L)" (Fermented Butter)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H71").Value = 775.75
$ws.Range("I71").Value = 524.9524
$ws.Range("J71").Value = 1202.7838
$ws.Range("K71").Value = 4724.5716
$ws.Range("L71").Value = 10825.0542
$ws.Range("M71").Value = -668.5716000000002
$ws.Range("N71").Value = -18937.0542

# GSM!130 - "Planisphere to Paper" (Chondrite Magitek Planisphere)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H130").Value = 30000
$ws.Range("J130").Value = 30000
$ws.Range("L130").Value = 30000
$ws.Range("N130").Value = -40040

# GSM!132 - "On Board for Lar" (Lar Ingot)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 4089.8
$ws.Range("I132").Value = 3289.1428
$ws.Range("J132").Value = 5958
$ws.Range("K132").Value = 9867.428400000001
$ws.Range("L132").Value = 17874
$ws.Range("M132").Value = -7337.428400000001
$ws.Range("N132").Value = -22934

# LTW!40 - "Best Served Toad" (Toad Leather)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 50004
$ws.Range("J40").Value = 0
$ws.Range("L40").Value = 0
$ws.Range("N40").ClearContents()

# LTW!122 - "Hell on Leather" (Gaja Leather)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 4189.722
$ws.Range("I122").Value = 4547.759
$ws.Range("K122").Value = 13643.277
$ws.Range("M122").Value = -11193.277

# WVR!47 - "The Wages of Sin" (Linen Coatee of Crafting)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H47").Value = 0
$ws.Range("J47").Value = 0
$ws.Range("L47").Value = 0
$ws.Range("N47").ClearContents()

# WVR!122 - "Heavy Armoire" (Dark Hempen Cloth)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1735.2941
$ws.Range("I122").Value = 1150
$ws.Range("J122").Value = 1813.3334
$ws.Range("K122").Value = 3450
$ws.Range("L122").Value = 5440.0002
$ws.Range("M122").Value = -1000
$ws.Range("N122").Value = -10340.0002

# WVR!124 - "Hot Heads" (Almasty Serge Hat of Casting)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H124").Value = 1689496.6
$ws.Range("J124").Value = 1689496.6
$ws.Range("L124").Value = 1689496.6
$ws.Range("N124").Value = -1699316.6
